$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws.Range("B9").Value = "Alvearie Team"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"
$ws.Rows.Item(11).Delete()
